{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// Summary of the change (per the diff):\n//   1. \"To do:\" splits into \"To \" / spell-checked \"do\" / \":\"\n//   2. In the big \"-tulelle oma damagemuoto...\" paragraph, several Finnish /\n//      compound words get wrapped with <w:proofErr> spell-check markers\n//      (damagemuoto, trigger-kolliderit, tuli+effekti, tallbushtest,\n//      damagesta, collideri, ei-triggeri) -- the run boundaries shift\n//      accordingly but the visible text is unchanged.\n//   3. The whole \"-kamerassa on verieffekti...\" paragraph is deleted.\n//   4. In the \"-laatikkojen...\" paragraph, \"sopis\" gets wrapped with\n//      proofErr marks.\n//   5. \"Bugeja:\" splits into spell-checked \"Bugeja\" / \":\"\n//\n// Because <w:proofErr/> is a presentation-only marker that isn't exposed as\n// a first class object in the Word JS API, we rebuild each affected\n// paragraph's contents from raw OOXML (Office.js's insertOoxml) so the\n// exact run/proofErr structure from the diff is reproduced, and we delete\n// the removed paragraph outright.\n\nfunction pkg(bodyInnerXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyInnerXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst PARA0_OOXML = pkg(\n  '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">To </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>do</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst PARA_DAMAGE_OOXML = pkg(\n  '<w:p>' +\n    '<w:r><w:t>-</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">tulelle oma </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>damagemuoto</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>, joka tulee r\\u00e4j\\u00e4hdyksist\\u00e4, tuliluodeista ja muista tulista</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> ja joka laskee my\\u00f6s </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>trigger-kolliderit</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> mukaan</w:t></w:r>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> Nyt </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>tuli</w:t></w:r>' +\n    '<w:r><w:t>effekti</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> toimii tuossa &quot;</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>tallbushtest</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>&quot;-puskassa, mutta se syt</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">tyy kaikesta </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>damagesta</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> ja </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>collideri</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> on </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>ei-triggeri</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>.</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst PARA_LAATIKKO_OOXML = pkg(\n  '<w:p>' +\n    '<w:r><w:t>-l</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">aatikkojen tuhoutumisessa </w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">ty\\u00f6nn\\u00f6n suunta </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>sopis</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> olla luodin mukainen.</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst PARA_BUGEJA_OOXML = pkg(\n  '<w:p>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Bugeja</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate each target paragraph by its (pre-edit) text so the script is\n// resilient to exact indices.\nlet pToDo = null;\nlet pDamage = null;\nlet pKamerassa = null;\nlet pLaatikko = null;\nlet pBugeja = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t === \"To do:\") {\n    pToDo = p;\n  } else if (t.indexOf(\"tulelle oma damagemuoto\") !== -1) {\n    pDamage = p;\n  } else if (t.indexOf(\"-kamerassa on\") !== -1) {\n    pKamerassa = p;\n  } else if (t.indexOf(\"laatikkojen tuhoutumisessa\") !== -1) {\n    pLaatikko = p;\n  } else if (t === \"Bugeja:\") {\n    pBugeja = p;\n  }\n}\n\nif (pKamerassa) {\n  pKamerassa.delete();\n}\nif (pToDo) {\n  pToDo.insertOoxml(PARA0_OOXML, Word.InsertLocation.replace);\n}\nif (pDamage) {\n  pDamage.insertOoxml(PARA_DAMAGE_OOXML, Word.InsertLocation.replace);\n}\nif (pLaatikko) {\n  pLaatikko.insertOoxml(PARA_LAATIKKO_OOXML, Word.InsertLocation.replace);\n}\nif (pBugeja) {\n  pBugeja.insertOoxml(PARA_BUGEJA_OOXML, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Summary of the change (per the diff):\n#   1. \"To do:\" splits into \"To \" / spell-checked \"do\" / \":\"\n#   2. In the big \"-tulelle oma damagemuoto...\" paragraph, several Finnish /\n#      compound words get wrapped with <w:proofErr> spell-check markers\n#      (damagemuoto, trigger-kolliderit, tuli+effekti, tallbushtest,\n#      damagesta, collideri, ei-triggeri) -- the run boundaries shift\n#      accordingly but the visible text is unchanged.\n#   3. The whole \"-kamerassa on verieffekti...\" paragraph is deleted.\n#   4. In the \"-laatikkojen...\" paragraph, \"sopis\" gets wrapped with\n#      proofErr marks.\n#   5. \"Bugeja:\" splits into spell-checked \"Bugeja\" / \":\"\n#\n# <w:proofErr/> is a presentation-only marker Word's proofer inserts; it is\n# not exposed as a first-class COM object, so each affected paragraph's\n# Range is rebuilt from raw OOXML via Range.InsertXML (which replaces the\n# exact Range it is called on), and the removed paragraph's Range is\n# deleted outright.\n\n$d = $word.ActiveDocument\n\nfunction New-PkgXml($bodyInnerXml) {\n    return '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyInnerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$PARA_TODO_XML = New-PkgXml('<w:p><w:r><w:t xml:space=\"preserve\">To </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>do</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>:</w:t></w:r></w:p>')\n\n$PARA_DAMAGE_XML = New-PkgXml('<w:p><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space=\"preserve\">tulelle oma </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>damagemuoto</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>, joka tulee r\u00e4j\u00e4hdyksist\u00e4, tuliluodeista ja muista tulista</w:t></w:r><w:r><w:t xml:space=\"preserve\"> ja joka laskee my\u00f6s </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>trigger-kolliderit</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> mukaan</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Nyt </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>tuli</w:t></w:r><w:r><w:t>effekti</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> toimii tuossa &quot;</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>tallbushtest</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>&quot;-puskassa, mutta se syt</w:t></w:r><w:r><w:t xml:space=\"preserve\">tyy kaikesta </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>damagesta</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> ja </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>collideri</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> on </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>ei-triggeri</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p>')\n\n$PARA_LAATIKKO_XML = New-PkgXml('<w:p><w:r><w:t>-l</w:t></w:r><w:r><w:t xml:space=\"preserve\">aatikkojen tuhoutumisessa </w:t></w:r><w:r><w:t xml:space=\"preserve\">ty\u00f6nn\u00f6n suunta </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>sopis</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> olla luodin mukainen.</w:t></w:r></w:p>')\n\n$PARA_BUGEJA_XML = New-PkgXml('<w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Bugeja</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>:</w:t></w:r></w:p>')\n\n# First pass: find and delete the \"-kamerassa on verieffekti...\" paragraph\n# (do deletions before other edits so paragraph indices of the remaining\n# targets are not disturbed by subsequent InsertXML calls).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*-kamerassa on*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# Second pass: rewrite the other affected paragraphs in place via InsertXML.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -eq \"To do:`r\") {\n        $p.Range.InsertXML($PARA_TODO_XML)\n    } elseif ($t -like \"*tulelle oma damagemuoto*\") {\n        $p.Range.InsertXML($PARA_DAMAGE_XML)\n    } elseif ($t -like \"*laatikkojen tuhoutumisessa*\") {\n        $p.Range.InsertXML($PARA_LAATIKKO_XML)\n    } elseif ($t -eq \"Bugeja:`r\") {\n        $p.Range.InsertXML($PARA_BUGEJA_XML)\n    }\n}\n"}
